# Update gh-pages to output generated at 456a3b4
# This bumps the "想去人数" (interested count) values in column F
# across the four worksheets, mirroring duplicated event rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1615
$ws.Range("F3").Value = 3375
$ws.Range("F6").Value = 2411
$ws.Range("F9").Value = 263
$ws.Range("F12").Value = 1128
$ws.Range("F13").Value = 482
$ws.Range("F17").Value = 5057
$ws.Range("F20").Value = 3721
$ws.Range("F23").Value = 4075
$ws.Range("F24").Value = 5425
$ws.Range("F25").Value = 132
$ws.Range("F27").Value = 590
$ws.Range("F28").Value = 3463
$ws.Range("F36").Value = 90
$ws.Range("F41").Value = 960
$ws.Range("F42").Value = 960
$ws.Range("F43").Value = 547
$ws.Range("F44").Value = 69
$ws.Range("F45").Value = 2542
$ws.Range("F47").Value = 202
$ws.Range("F49").Value = 3790

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F22").Value = 51

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3035

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3035
$ws.Range("F3").Value = 1615
$ws.Range("F4").Value = 3375
$ws.Range("F7").Value = 2411
$ws.Range("F10").Value = 263
$ws.Range("F14").Value = 1128
$ws.Range("F15").Value = 482
$ws.Range("F19").Value = 5057
$ws.Range("F21").Value = 4075
$ws.Range("F22").Value = 5425
$ws.Range("F23").Value = 132
$ws.Range("F25").Value = 590
$ws.Range("F26").Value = 3463
$ws.Range("F34").Value = 90
$ws.Range("F39").Value = 960
$ws.Range("F41").Value = 547
$ws.Range("F43").Value = 69
$ws.Range("F44").Value = 51
$ws.Range("F45").Value = 2542
$ws.Range("F47").Value = 202
$ws.Range("F49").Value = 3790
